# Apply updates described by the diff to 杭州-漫展信息.xlsx
#
# Sheet "展览" (sheet1) and "全部类型" (sheet4) both list the same events.
# Most changes are simple "想去人数" (want-to-go count) bumps in column F.
# In addition, sheet "全部类型" is missing a couple of row updates that
# "展览" already reflects: an old, already-finished event row needs to be
# replaced by newer event data (equivalent to deleting the stale row and
# shifting everything up), a brand new ticket row needs to be added, and a
# stale row at the bottom needs to be replaced as well. Because column A
# holds a hard-coded sequential index (not a formula), we implement all of
# this purely by overwriting cell contents row-by-row rather than using
# real Insert/Delete row operations (which would shift the hard-coded
# index values along with the rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "展览" -- simple "想去人数" (column F) bumps
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")

$exhibitUpdates = @{
    8  = 640
    13 = 571
    17 = 1775
    18 = 1456
    22 = 312
    23 = 519
    28 = 2626
    33 = 27
    40 = 647
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F" + $row).Value = $exhibitUpdates[$row]
}

# ---------------------------------------------------------------------
# 2) Sheet "全部类型" -- simple "想去人数" (column F) bumps
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allTypesUpdates = @{
    9  = 640
    18 = 571
    23 = 1775
    24 = 1456
    27 = 312
    29 = 519
    37 = 27
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F" + $row).Value = $allTypesUpdates[$row]
}

# ---------------------------------------------------------------------
# Helper to write a full data row (columns B..I) on sheet "全部类型".
# Column B holds plain text dates such as "2024-03-23"; writing that
# directly through .Value would make Excel auto-convert it into a date
# serial number, so column B is explicitly pre-formatted as Text first.
#
# NOTE: this runtime's PowerShell engine does not bind named
# (-Param value) arguments to function "param()" blocks correctly, so
# this helper is called using plain positional arguments only.
# ---------------------------------------------------------------------
function Set-AllTypesRow($Row, $B, $C, $D, $E, $F, $G, $H, $I) {
    $bCell = $wsAll.Range("B" + $Row)
    $bCell.NumberFormat = "@"
    $bCell.Value = $B

    $wsAll.Range("C" + $Row).Value = $C
    $wsAll.Range("D" + $Row).Value = $D
    $wsAll.Range("E" + $Row).Value = $E
    $wsAll.Range("F" + $Row).Value = $F
    $wsAll.Range("G" + $Row).Value = $G
    $wsAll.Range("H" + $Row).Value = $H
    $wsAll.Range("I" + $Row).Value = $I
}

# ---------------------------------------------------------------------
# 3) Sheet "全部类型" rows 32-35: the stale "2024-03-17 ComicMe · 马正阳
#    专场" row (already passed) is removed, which shifts the AD02 rows
#    up by one; a brand-new "亦之紫F、L句号内场票" row is introduced
#    before the existing "钟晨瑶内场票" row (row 36, which keeps its
#    current content and is intentionally left untouched below).
# ---------------------------------------------------------------------

Set-AllTypesRow 32 `
    "2024-03-23" `
    "杭州·AD02动漫展" `
    "浙江省杭州市萧山区奔竞大道353号 国际博览中心" `
    "2024.03.23 10:00-03.24 17:00" `
    2626 `
    75 `
    "https://show.bilibili.com/platform/detail.html?id=80905" `
    "//i1.hdslb.com/bfs/openplatform/202401/D3QaPamg1705397424553.jpeg"

Set-AllTypesRow 33 `
    "2024-03-23" `
    "杭州·AD02动漫展  青柳尊哉内场票" `
    "浙江省杭州市萧山区奔竞大道353号 国际博览中心" `
    "2024.03.23 10:00-03.23 17:00" `
    173 `
    528 `
    "https://show.bilibili.com/platform/detail.html?id=81503" `
    "//i1.hdslb.com/bfs/openplatform/202401/OmqxboDC1706522627528.jpeg"

Set-AllTypesRow 34 `
    "2024-03-24" `
    "杭州·AD02动漫展  岩永彻也内场票" `
    "浙江省杭州市萧山区奔竞大道353号 国际博览中心" `
    "2024.03.24 10:00-03.24 17:00" `
    99 `
    528 `
    "https://show.bilibili.com/platform/detail.html?id=81239" `
    "//i0.hdslb.com/bfs/openplatform/202401/hww9WUpD1705914756383.jpeg"

Set-AllTypesRow 35 `
    "2024-03-24" `
    "杭州·AD02动漫展--亦之紫F、L句号内场票" `
    "钱江世纪城奔竞大道353号 杭州国际博览中心" `
    "2024.03.24 12:00-03.24 16:00" `
    55 `
    258 `
    "https://show.bilibili.com/platform/detail.html?id=81836" `
    "//i1.hdslb.com/bfs/openplatform/202402/ecrRfQce1707375167618.jpeg"

# Row 36 ("钟晨瑶内场票") is unchanged by the diff, so it is left as-is.

# ---------------------------------------------------------------------
# 4) Sheet "全部类型" rows 43-44: a new "倒霉死勒内场票" row is inserted
#    before the existing "钱琛签售礼包" row (which shifts from row 43 to
#    44), and the stale "赛马娘only" row (old row 44, already passed) is
#    removed / overwritten in the process.
# ---------------------------------------------------------------------

Set-AllTypesRow 43 `
    "2024-04-05" `
    "杭州·ESCC电竞博览会 倒霉死勒内场票" `
    "钱江世纪城奔竞大道353号 杭州国际博览中心" `
    "2024.04.05 09:30-04.05 16:30" `
    647 `
    248 `
    "https://show.bilibili.com/platform/detail.html?id=81681" `
    "//i1.hdslb.com/bfs/openplatform/202402/suXI547M1706862164353.png"

Set-AllTypesRow 44 `
    "2024-04-05" `
    "杭州·ESCC电竞博览会·钱琛签售礼包" `
    "钱江世纪城奔竞大道353号 杭州国际博览中心" `
    "2024.04.05 09:30-04.05 16:30" `
    83 `
    39.9 `
    "https://show.bilibili.com/platform/detail.html?id=81680" `
    "//i2.hdslb.com/bfs/openplatform/202402/noqtqw701706861615316.png"

Write-Output "Edit applied"
